# Add a new "Madness Domain" sub-class entry to the Cleric > Domain list.
# This inserts a new row above the existing row 33 ("Nature Domain"),
# shifting every row below it down by one, and fills the new B33 cell
# with "Madness Domain" (alphabetically between "Light Domain" and
# "Nature Domain"). Finally, move the active selection to B40 (mirroring
# the author's on-screen selection after the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(33).Insert()
$ws.Cells.Item(33, 2).Value = "Madness Domain"

$ws.Range("B40").Select()
